$wb = $excel.ActiveWorkbook

# Rename existing sheets first (order: DOANH SỐ CÁ NHÂN -> DAONH SỐ CÁ NHÂN)
$personal = $wb.Worksheets.Item("DOANH SỐ CÁ NHÂN")
$personal.Name = "DAONH SỐ CÁ NHÂN"

# Insert a brand-new sheet at the very front of the workbook for the
# detailed revenue report "CHI TIẾT DOANH THU".
$detail = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$detail.Name = "CHI TIẾT DOANH THU"

# Header row
$headers = @("Tiền tố", "Mã dịch vụ", "Ngày thực hiện", "Cơ sở", "Tên dịch vụ", "Khách hàng", "Nguồn khách", "Sale chính", "Đơn giá gốc", "Sale phụ", "Upsale", "Đơn giá", "Bác sĩ 1", "Bác sĩ 2", "Thanh toán lần đầu", "Trả sau", "Đã thanh toán", "Dư nợ", "Phụ phẫu 1", "Phụ phẫu 2")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $detail.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2 (date column is kept as literal text, like the source data, via a
# leading apostrophe so Excel doesn't auto-convert it to a date serial;
# blank numeric cells in the source report are stored as 0)
$row2 = @("HD-LUXURY", 507, "'07-01-2024", "LONG XUYÊN", "Tiềm cằm", "Nana", "Khách cũ", "Nguyễn Phúc Nam", 3000000, 0, 0, 3000000, "Đặng Ngọc Mai", 0, 3000000, 0, 3000000, 0, "Đào Vương Anh", 0)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $detail.Cells.Item(2, $i + 1).Value = $row2[$i]
}

# Row 3
$row3 = @("HD-LUXURY", 511, "'07-02-2024", "LONG XUYÊN", "Điêu khắc mày", "Nguyễn ngọc trinh", "CTV", 0, 1500000, 0, 0, 1500000, "Đặng Ngọc Mai", 0, 1500000, 0, 1500000, 0, 0, 0)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $detail.Cells.Item(3, $i + 1).Value = $row3[$i]
}

$detail.Range("A1").Select()
